$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 150.46666
$ws.Range("I6").Value = 125.35714
$ws.Range("K6").Value = 376.07142
$ws.Range("M6").Value = -264.07142
$ws.Range("H31").Value = 4819.3335
$ws.Range("I31").Value = 2979
$ws.Range("K31").Value = 8937
$ws.Range("M31").Value = -8707
$ws.Range("H41").Value = 1873
$ws.Range("I41").Value = 2941.75
$ws.Range("J41").Value = 923
$ws.Range("K41").Value = 2941.75
$ws.Range("L41").Value = 923
$ws.Range("M41").Value = -2501.75
$ws.Range("N41").Value = -1803
$ws.Range("H88").Value = 1063.5
$ws.Range("I88").Value = 1208.3334
$ws.Range("J88").Value = 846.25
$ws.Range("K88").Value = 1208.3334
$ws.Range("L88").Value = 846.25
$ws.Range("M88").Value = -802.3334
$ws.Range("N88").Value = -1658.25
$ws.Range("H91").Value = 1063.5
$ws.Range("I91").Value = 1208.3334
$ws.Range("J91").Value = 846.25
$ws.Range("K91").Value = 1208.3334
$ws.Range("L91").Value = 846.25
$ws.Range("M91").Value = 195.6666
$ws.Range("N91").Value = -3654.25
$ws.Range("H101").Value = 0
$ws.Range("I101").Value = 0
$ws.Range("K101").Value = 0
$ws.Range("M101").ClearContents()
$ws.Range("H116").Value = 19000
$ws.Range("J116").Value = 0
$ws.Range("L116").Value = 0
$ws.Range("N116").ClearContents()
$ws.Range("H132").Value = 4435.8184
$ws.Range("I132").Value = 4809.4
$ws.Range("J132").Value = 700
$ws.Range("K132").Value = 14428.2
$ws.Range("L132").Value = 2100
$ws.Range("M132").Value = -11898.2
$ws.Range("N132").Value = -7160
$ws.Range("H137").Value = 3374.8572
$ws.Range("I137").Value = 3325
$ws.Range("J137").Value = 3499.5
$ws.Range("K137").Value = 9975
$ws.Range("L137").Value = 10498.5
$ws.Range("M137").Value = -7425
$ws.Range("N137").Value = -15598.5
$ws.Range("H138").Value = 2445.8635
$ws.Range("I138").Value = 968.3333
$ws.Range("J138").Value = 2999.9375
$ws.Range("K138").Value = 2904.9999
$ws.Range("L138").Value = 8999.8125
$ws.Range("M138").Value = 2235.0001
$ws.Range("N138").Value = -19279.8125
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 1236.625
$ws.Range("I74").Value = 1236.625
$ws.Range("K74").Value = 1236.625
$ws.Range("M74").Value = -362.625
$ws.Range("H77").Value = 1236.625
$ws.Range("I77").Value = 1236.625
$ws.Range("K77").Value = 6183.125
$ws.Range("M77").Value = -1815.125
$ws.Range("H96").Value = 33137.8
$ws.Range("J96").Value = 33137.8
$ws.Range("L96").Value = 33137.8
$ws.Range("N96").Value = -38629.8
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H26").Value = 15499.75
$ws.Range("I26").Value = 15499.75
$ws.Range("J26").Value = 0
$ws.Range("K26").Value = 15499.75
$ws.Range("L26").Value = 0
$ws.Range("M26").Value = -15207.75
$ws.Range("N26").ClearContents()
$ws.Range("H31").Value = 0
$ws.Range("I31").Value = 0
$ws.Range("K31").Value = 0
$ws.Range("M31").ClearContents()
$ws.Range("H106").Value = 35000
$ws.Range("J106").Value = 35000
$ws.Range("L106").Value = 35000
$ws.Range("N106").Value = -37524
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1814.8823
$ws.Range("J31").Value = 2499.4
$ws.Range("L31").Value = 2499.4
$ws.Range("N31").Value = -3089.4
$ws.Range("H34").Value = 1814.8823
$ws.Range("J34").Value = 2499.4
$ws.Range("L34").Value = 2499.4
$ws.Range("N34").Value = -2903.4
$ws.Range("H59").Value = 60000
$ws.Range("J59").Value = 60000
$ws.Range("L59").Value = 60000
$ws.Range("N59").Value = -62290
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H109").Value = 0
$ws.Range("I109").Value = 0
$ws.Range("K109").Value = 0
$ws.Range("M109").ClearContents()
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H49").Value = 5000
$ws.Range("I49").Value = 2000
$ws.Range("J49").Value = 8000
$ws.Range("K49").Value = 2000
$ws.Range("L49").Value = 8000
$ws.Range("M49").Value = -1816
$ws.Range("N49").Value = -8368
$ws.Range("H95").Value = 40488.75
$ws.Range("J95").Value = 40488.75
$ws.Range("L95").Value = 40488.75
$ws.Range("N95").Value = -45980.75
$ws.Range("H104").Value = 0
$ws.Range("J104").Value = 0
$ws.Range("L104").Value = 0
$ws.Range("N104").ClearContents()
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 1409.8
$ws.Range("I7").Value = 1409.8
$ws.Range("K7").Value = 1409.8
$ws.Range("M7").Value = -1297.8
$ws.Range("H16").Value = 2557.8
$ws.Range("I16").Value = 2933
$ws.Range("J16").Value = 1995
$ws.Range("K16").Value = 2933
$ws.Range("L16").Value = 1995
$ws.Range("M16").Value = -2763
$ws.Range("N16").Value = -2335
$ws.Range("H25").Value = 3008
$ws.Range("I25").Value = 0
$ws.Range("J25").Value = 3008
$ws.Range("K25").Value = 0
$ws.Range("L25").Value = 3008
$ws.Range("M25").ClearContents()
$ws.Range("N25").Value = -3468
$ws.Range("H46").Value = 4289.6
$ws.Range("I46").Value = 3483
$ws.Range("J46").Value = 5499.5
$ws.Range("K46").Value = 3483
$ws.Range("L46").Value = 5499.5
$ws.Range("M46").Value = -3295
$ws.Range("N46").Value = -5875.5
$ws.Range("H122").Value = 7500
$ws.Range("I122").Value = 0
$ws.Range("J122").Value = 7500
$ws.Range("K122").Value = 0
$ws.Range("L122").Value = 22500
$ws.Range("M122").ClearContents()
$ws.Range("N122").Value = -27400
$ws.Range("H126").Value = 1409.8
$ws.Range("I126").Value = 1409.8
$ws.Range("K126").Value = 4229.4
$ws.Range("M126").Value = -1759.4
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H3").Value = 1502001
$ws.Range("I3").Value = 1502001
$ws.Range("K3").Value = 1502001
$ws.Range("M3").Value = -1501887
$ws.Range("H9").Value = 10000
$ws.Range("I9").Value = 0
$ws.Range("J9").Value = 10000
$ws.Range("K9").Value = 0
$ws.Range("L9").Value = 10000
$ws.Range("M9").ClearContents()
$ws.Range("N9").Value = -10280
$ws.Range("H28").Value = 149999.5
$ws.Range("J28").Value = 150000
$ws.Range("L28").Value = 150000
$ws.Range("N28").Value = -150696
$ws.Range("H105").Value = 67449.5
$ws.Range("J105").Value = 67449.5
$ws.Range("L105").Value = 67449.5
$ws.Range("N105").Value = -74437.5
$ws.Range("H122").Value = 0
$ws.Range("I122").Value = 0
$ws.Range("K122").Value = 0
$ws.Range("M122").ClearContents()
$ws.Range("H126").Value = 11000
$ws.Range("J126").Value = 10000
$ws.Range("L126").Value = 30000
$ws.Range("N126").Value = -34940
$ws.Range("H136").Value = 1396.5834
$ws.Range("I136").Value = 1069
$ws.Range("K136").Value = 3207
$ws.Range("M136").Value = -657
